# Gym data workbook update:
#  - Refresh a handful of placeholder/outdated prices on existing sheets
#  - Add a new "Leg Extensions" sheet with product data + hyperlinks

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update existing price cells
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
$ws1.Range("C2").Value = '''$2,072.00'
$ws1.Range("C3").Value = '''$1,119.99'

$ws2 = $wb.Worksheets.Item("Squat Stands")
$ws2.Range("C2").Value = '''$1,498.00'
$ws2.Range("C3").Value = '''$459.99'

# ---------------------------------------------------------------------
# 2) Add the new "Leg Extensions" sheet after "Squat Stands"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Leg Extensions"

# Column widths (matches the other two sheets)
$ws3.Columns.Item(1).ColumnWidth = 25
$ws3.Columns.Item(2).ColumnWidth = 20
$ws3.Columns.Item(3).ColumnWidth = 15
$ws3.Columns.Item(4).ColumnWidth = 12
$ws3.Columns.Item(5).ColumnWidth = 25
$ws3.Columns.Item(6).ColumnWidth = 35

# Header row
$ws3.Range("A1").Value = "Product Name"
$ws3.Range("B1").Value = "Manufacturer"
$ws3.Range("C1").Value = "Price"
$ws3.Range("D1").Value = "Country"
$ws3.Range("E1").Value = "Image URL"
$ws3.Range("F1").Value = "Product Page"
$ws3.Range("A1:F1").Font.Bold = $true
$ws3.Range("A1:F1").HorizontalAlignment = -4108
$ws3.Range("A1:F1").Borders.LineStyle = 1

# Row 2 - Oak Club MFG (no product available)
$ws3.Range("A2").Value = "No Leg Extension Available"
$ws3.Range("B2").Value = "Oak Club MFG"
$ws3.Range("C2").Value = "Not Available"
$ws3.Range("D2").Value = "Canada"
$ws3.Range("E2").Value = "https://t3.ftcdn.net/jpg/01/12/43/90/360_F_112439022_Sft6cXK9GLnzWjjIkVMj2Lt34RcKUpxm.jpg"
$ws3.Range("F2").Value = "Not Available"

# Row 3 - Titan Fitness
$ws3.Range("A3").Value = "Leg Extension and Curl Machine | 10 - 250 LB Selector"
$ws3.Range("B3").Value = "Titan Fitness"
$ws3.Range("C3").Value = '''$2,529.99'
$ws3.Range("D3").Value = "China"
$ws3.Range("E3").Value = "https://titan.fitness/cdn/shop/files/401926_01.jpg?v=1739302160&width=832"
$ws3.Range("F3").Value = "https://titan.fitness/products/selectorized-leg-extension-and-curl-machine"

# Row 4 - Rogue Fitness (no product available)
$ws3.Range("A4").Value = "No Leg Extension Available"
$ws3.Range("B4").Value = "Rogue Fitness"
$ws3.Range("C4").Value = "Not Available"
$ws3.Range("D4").Value = "USA"
$ws3.Range("E4").Value = "https://t3.ftcdn.net/jpg/01/12/43/90/360_F_112439022_Sft6cXK9GLnzWjjIkVMj2Lt34RcKUpxm.jpg"
$ws3.Range("F4").Value = "Not Available"

# Row 5 - Stray Dog Strength
$ws3.Range("A5").Value = "Selectorized Seated Leg Curl/Extension"
$ws3.Range("B5").Value = "Stray Dog Strength"
$ws3.Range("C5").Value = '''$5,250.00'
$ws3.Range("D5").Value = "USA"
$ws3.Range("E5").Value = "https://shop.straydogstrength.com/cdn/shop/files/2325-RIGHT-RED_eee5d4da-9504-4bb9-b7e3-f98e7e85c231.jpg?v=1743705611&width=823"
$ws3.Range("F5").Value = "https://shop.straydogstrength.com/products/selectorized-seated-leg-curl-extension"

# Row 6 - Sorinex
$ws3.Range("A6").Value = "LEG EXTENSION / CURL MACHINE"
$ws3.Range("B6").Value = "Sorinex"
$ws3.Range("C6").Value = '''$4,149.00'
$ws3.Range("D6").Value = "USA"
$ws3.Range("E6").Value = "https://cdn.shopify.com/s/files/1/2559/4942/files/LegCurlLegExtension.jpg?v=1733930466"
$ws3.Range("F6").Value = "https://www.sorinex.com/products/leg-extension-curl-machine?Title=Default+Title"

# Hyperlinks (image URL + product page columns, matching other sheets)
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://t3.ftcdn.net/jpg/01/12/43/90/360_F_112439022_Sft6cXK9GLnzWjjIkVMj2Lt34RcKUpxm.jpg")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://titan.fitness/cdn/shop/files/401926_01.jpg?v=1739302160&width=832")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://titan.fitness/products/selectorized-leg-extension-and-curl-machine")
$ws3.Hyperlinks.Add($ws3.Range("E4"), "https://t3.ftcdn.net/jpg/01/12/43/90/360_F_112439022_Sft6cXK9GLnzWjjIkVMj2Lt34RcKUpxm.jpg")
$ws3.Hyperlinks.Add($ws3.Range("E5"), "https://shop.straydogstrength.com/cdn/shop/files/2325-RIGHT-RED_eee5d4da-9504-4bb9-b7e3-f98e7e85c231.jpg?v=1743705611&width=823")
$ws3.Hyperlinks.Add($ws3.Range("F5"), "https://shop.straydogstrength.com/products/selectorized-seated-leg-curl-extension")
$ws3.Hyperlinks.Add($ws3.Range("E6"), "https://cdn.shopify.com/s/files/1/2559/4942/files/LegCurlLegExtension.jpg?v=1733930466")
$ws3.Hyperlinks.Add($ws3.Range("F6"), "https://www.sorinex.com/products/leg-extension-curl-machine?Title=Default+Title")

$ws3.Range("A1").Select()
